# Apply the "Saldo_guide" daily refresh update:
#  - sheet re-named to reflect the new export timestamp (20240807 -> 20240808)
#  - every row's "Dt. Referencia" (column G) rolls from 45511 to 45512 (one day later)
#  - two accounts (row 17 and row 112) had their balance refreshed, so
#    "Saldo Previsto" (E) and "Vl. Total" (H) change together on those rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new export run.
$ws.Name = "IClientBalance-20240808-103649-"

# Column G holds the reference date as an Excel serial number; bump every
# data row (2 through 274) from 45511 (2024-08-07) to 45512 (2024-08-08).
$lastRow = 274
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45512
}

# Row 17: balance updated from 708.26 to 680.86 (columns E and H mirror each other).
$ws.Cells.Item(17, 5).Value = 680.86
$ws.Cells.Item(17, 8).Value = 680.86

# Row 112: balance updated from 34.15 to 459.25 (columns E and H mirror each other).
$ws.Cells.Item(112, 5).Value = 459.25
$ws.Cells.Item(112, 8).Value = 459.25
